$p = $ppt.ActivePresentation
$s1 = $p.Slides.Item(1)
$tcs = $s1.ThemeColorScheme

function Set-RGB($idx, $hex) {
  $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
  $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
  $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
  $val = $r + ($g * 256) + ($b * 65536)
  $tcs.Colors($idx).RGB = $val
}

Set-RGB 1 "000000"
Set-RGB 2 "FFFFFF"
Set-RGB 3 "44546A"
Set-RGB 4 "E7E6E6"
Set-RGB 5 "5B9BD5"
Set-RGB 6 "ED7D31"
Set-RGB 7 "A5A5A5"
Set-RGB 8 "FFC000"
Set-RGB 9 "4472C4"
Set-RGB 10 "70AD47"
Set-RGB 11 "0563C1"
Set-RGB 12 "954F72"

Write-Output "done"
